# Week 17 data log for the Packers Players Data workbook.
#
# Sheet "Rushing" (sheet1): a new row for J.Love is inserted right after
# A.Rodgers. Rather than using Rows.Insert() (which can invent new,
# unreferenced cell styles), we shift every player's rushing stats down
# one row by literal value assignment, fill in J.Love's (week-17-only)
# line, and append D.Dafney's updated totals as a brand-new last row,
# copying the existing index-column style onto the newly created cells.
#
# Sheet "Receiving" (sheet2): same 14 players, just refreshed
# cumulative-through-week-17 totals -- no row shuffling needed there.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Rushing")
$ws2 = $wb.Worksheets.Item("Receiving")

# ---------------------------------------------------------------------
# Rushing sheet
# ---------------------------------------------------------------------

# A.Rodgers (row 2) picked up more work this week.
$ws1.Range("D2").Value = 7
$ws1.Range("E2").Value = 8
$ws1.Range("F2").Value = 8

# Every other rushing row shifts down by one to make room for J.Love,
# who slots in right under A.Rodgers. Values below are each player's
# updated (through week 17) rushing line.
$ws1.Range("B3").Value = "J.Love"
$ws1.Range("C3").Value = 0
$ws1.Range("D3").Value = 0
$ws1.Range("E3").Value = 1
$ws1.Range("F3").Value = 0

$ws1.Range("B4").Value = "A.Jones"
$ws1.Range("C4").Value = 112
$ws1.Range("D4").Value = 55
$ws1.Range("E4").Value = 9
$ws1.Range("F4").Value = 34

$ws1.Range("B5").Value = "A.Dillon"
$ws1.Range("C5").Value = 86
$ws1.Range("D5").Value = 62
$ws1.Range("E5").Value = 16
$ws1.Range("F5").Value = 36

$ws1.Range("B6").Value = "K.Hill"
$ws1.Range("C6").Value = 3
$ws1.Range("D6").Value = 5
$ws1.Range("E6").Value = 0
$ws1.Range("F6").Value = 0

$ws1.Range("B7").Value = "P.Taylor"
$ws1.Range("C7").Value = 5
$ws1.Range("D7").Value = 3
$ws1.Range("E7").Value = 0
$ws1.Range("F7").Value = 2

$ws1.Range("B8").Value = "A.Lazard"
$ws1.Range("C8").Value = 2
$ws1.Range("D8").Value = 1
$ws1.Range("E8").Value = 0
$ws1.Range("F8").Value = 0

$ws1.Range("B9").Value = "Am.Rodgers"
$ws1.Range("C9").Value = 1
$ws1.Range("D9").Value = 0
$ws1.Range("E9").Value = 0
$ws1.Range("F9").Value = 0

$ws1.Range("B10").Value = "E.St. Brown"
$ws1.Range("C10").Value = 3
$ws1.Range("D10").Value = 0
$ws1.Range("E10").Value = 0
$ws1.Range("F10").Value = 0

$ws1.Range("B11").Value = "J.Deguara"
$ws1.Range("C11").Value = 1
$ws1.Range("D11").Value = 0
$ws1.Range("E11").Value = 0
$ws1.Range("F11").Value = 0

# D.Dafney's updated line lands on a brand-new row 12. Give its index
# cell (A12) the same look as the rest of the index column (A2:A11)
# by copying formatting from the cell directly above it, then fill in
# the values.
$ws1.Range("A11").Copy()
$ws1.Range("A12").PasteSpecial(-4122)
$ws1.Range("A12").Value = 10
$ws1.Range("B12").Value = "D.Dafney"
$ws1.Range("C12").Value = 1
$ws1.Range("D12").Value = 0
$ws1.Range("E12").Value = 0
$ws1.Range("F12").Value = 1

# ---------------------------------------------------------------------
# Receiving sheet (no new rows, just updated week-17 cumulative totals)
# ---------------------------------------------------------------------

$ws2.Range("C2").Value = 64
$ws2.Range("D2").Value = 53
$ws2.Range("G2").Value = 15
$ws2.Range("H2").Value = 12

$ws2.Range("C3").Value = 29
$ws2.Range("D3").Value = 24

$ws2.Range("C5").Value = 133
$ws2.Range("D5").Value = 110
$ws2.Range("E5").Value = 40
$ws2.Range("F5").Value = 30
$ws2.Range("G5").Value = 29
$ws2.Range("H5").Value = 24

$ws2.Range("C6").Value = 21
$ws2.Range("D6").Value = 16
$ws2.Range("E6").Value = 23
$ws2.Range("G6").Value = 5

$ws2.Range("C7").Value = 50
$ws2.Range("D7").Value = 39
$ws2.Range("E7").Value = 12
$ws2.Range("F7").Value = 8
$ws2.Range("G7").Value = 14
$ws2.Range("H7").Value = 7

$ws2.Range("C10").Value = 4
$ws2.Range("D10").Value = 3
$ws2.Range("G10").Value = 2

$ws2.Range("C11").Value = 9
$ws2.Range("D11").Value = 5

$ws2.Range("C12").Value = 3
$ws2.Range("D12").Value = 3
$ws2.Range("G12").Value = 1
$ws2.Range("H12").Value = 1

$ws2.Range("C14").Value = 24
$ws2.Range("D14").Value = 19
$ws2.Range("E14").Value = 3
$ws2.Range("G14").Value = 6
$ws2.Range("H14").Value = 3

$ws2.Range("C15").Value = 2
